# feat: add 2022-Q3 data
#
# This script:
#  1. Creates a new worksheet "2022-Q3" (as a copy of "2022-Q2", so that it
#     keeps exactly the same layout/styling) positioned right after "总计"
#     and before the existing "2022-Q2" sheet.
#  2. Overwrites the fund figures on that new sheet with the 2022-Q3 values.
#  3. Updates the "总计" (summary) sheet, inserting a new top data row for
#     2022-Q3 and shifting the existing 2022-Q2 / 2022-Q1 rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate "2022-Q2" -> new sheet placed right before it, then rename
# ---------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2) | Out-Null

$wsQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$wsQ3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2) Update the fund data on the new "2022-Q3" sheet
# ---------------------------------------------------------------------

# Columns D:G hold numbers formatted as plain text (e.g. "9.90"); force
# text storage (without leaving a visible "@" format behind) so the
# values keep being written as strings and not auto-converted to numbers.
$textRange = $wsQ3.Range("D2:G7")
$textRange.NumberFormat = "@"

$wsQ3.Range("D2").Value = "9.90"
$wsQ3.Range("E2").Value = "61.86"
$wsQ3.Range("F2").Value = "3.22"
$wsQ3.Range("G2").Value = "0.3188"
$wsQ3.Range("H2").Value = 5

$wsQ3.Range("D3").Value = "6.88"
$wsQ3.Range("E3").Value = "61.00"
$wsQ3.Range("F3").Value = "2.24"
$wsQ3.Range("G3").Value = "0.1541"
$wsQ3.Range("H3").Value = 10

$wsQ3.Range("D4").Value = "3.65"
$wsQ3.Range("E4").Value = "60.88"
$wsQ3.Range("F4").Value = "3.06"
$wsQ3.Range("G4").Value = "0.1117"
$wsQ3.Range("H4").Value = 8

$wsQ3.Range("C5").Value = "大成投资严选六月持有混合A"
$wsQ3.Range("D5").Value = "3.10"
$wsQ3.Range("E5").Value = "66.75"
$wsQ3.Range("F5").Value = "2.95"
$wsQ3.Range("G5").Value = "0.0914"
$wsQ3.Range("H5").Value = 7

$wsQ3.Range("C6").Value = "大成投资严选六月持有混合C"
$wsQ3.Range("D6").Value = "0.22"
$wsQ3.Range("E6").Value = "66.75"
$wsQ3.Range("F6").Value = "2.95"
$wsQ3.Range("G6").Value = "0.0065"
$wsQ3.Range("H6").Value = 7

$wsQ3.Range("D7").Value = "0.17"
$wsQ3.Range("E7").Value = "60.88"
$wsQ3.Range("F7").Value = "3.06"
$wsQ3.Range("G7").Value = "0.0052"
$wsQ3.Range("H7").Value = 8

# Drop the temporary "@" text format again so the cells end up unstyled,
# exactly like the rest of the sheet.
$textRange.Style = "Normal"

# ---------------------------------------------------------------------
# 3) Update the "总计" summary sheet: push the existing rows down one and
#    insert the new 2022-Q3 row at the top of the data.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Copy the formatting of row 3 down into row 4 for the row that is about
# to hold the (previously row 3) 2022-Q1 data.
$wsTotal.Range("A3:D3").Copy() | Out-Null
$wsTotal.Range("A4:D4").PasteSpecial(-4122) | Out-Null

# Row 4 becomes the old 2022-Q1 row (previously row 3)
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q1"
$wsTotal.Range("C4").Value = 4
$wsTotal.Range("D4").Value = 0.65

# Row 3 becomes the old 2022-Q2 row (previously row 2)
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 6
$wsTotal.Range("D3").Value = 0.8

# Row 2 becomes the new 2022-Q3 row
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 6
$wsTotal.Range("D2").Value = 0.69

# Restore "总计" as the active/selected sheet (as in the original workbook).
$wsTotal.Select()
